$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> hashtable of column letter -> new value, derived from the
# commit's updated cryptocurrency price / 1h volume figures.
$updates = @{
    2 = @{ D = "64.407.66"; E = "  -3.13%  " }
    3 = @{ D = "3.174.51"; E = "  -4.64%  " }
    4 = @{ E = "  +0.03%  " }
    5 = @{ D = "569.62"; E = "  -2.91%  " }
    6 = @{ D = "169.05"; E = "  -7.95%  " }
    7 = @{ D = "0.610"; E = "  -5.55%  " }
    8 = @{ E = "  +0.02%  " }
    9 = @{ D = "3.176.90"; E = "  -4.53%  " }
    10 = @{ D = "0.121"; E = "  -4.98%  " }
    11 = @{ D = "6.78"; E = "  +0.18%  " }
    13 = @{ D = "3.728.70"; E = "  -4.54%  " }
    14 = @{ E = "  -2.09%  " }
    15 = @{ D = "64.452.40"; E = "  -3.07%  " }
    16 = @{ D = "25.37"; E = "  -3.95%  " }
    17 = @{ E = "  -3.64%  " }
    18 = @{ D = "3.157.08"; E = "  -4.72%  " }
    19 = @{ D = "420.08"; E = "  -2.75%  " }
    20 = @{ D = "12.86"; E = "  -3.54%  " }
    21 = @{ D = "5.37"; E = "  -3.17%  " }
    22 = @{ D = "7.06"; E = "  -5.23%  " }
    23 = @{ E = "  -0.08%  " }
    24 = @{ D = "70.24"; E = "  -2.88%  " }
    25 = @{ D = "0.204"; E = "  +2.69%  " }
    26 = @{ D = "0.488"; E = "  -5.78%  " }
    27 = @{ D = "0.0000106"; E = "  -8.03%  " }
    28 = @{ D = "8.77"; E = "  -3.05%  " }
    29 = @{ E = "  -0.03%  " }
    30 = @{ E = "  -6.31%  " }
    31 = @{ D = "21.75"; E = "  -3.02%  " }
    33 = @{ D = "5.04"; E = "  -3.71%  " }
    34 = @{ D = "6.33"; E = "  -4.72%  " }
    35 = @{ E = "  -5.68%  " }
    36 = @{ D = "157.50"; E = "  -1.50%  " }
    37 = @{ D = "1.37"; E = "  -6.90%  " }
    38 = @{ D = "2.719.30"; E = "  -6.07%  " }
    39 = @{ E = "  -7.25%  " }
    40 = @{ D = "24.28"; E = "  -9.13%  " }
    41 = @{ E = "  -4.09%  " }
    42 = @{ D = "39.13"; E = "  -3.01%  " }
    43 = @{ D = "0.710"; E = "  -7.59%  " }
    44 = @{ D = "0.0621"; E = "  -7.12%  " }
    45 = @{ D = "5.60"; E = "  -6.81%  " }
    46 = @{ E = "  -4.30%  " }
    47 = @{ D = "21.69"; E = "  -7.82%  " }
    48 = @{ D = "292.46"; E = "  -8.07%  " }
    49 = @{ E = "  +0.04%  " }
    50 = @{ E = "  -13.79%  " }
    51 = @{ D = "0.0990"; E = "  -5.80%  " }
}

foreach ($rowNum in $updates.Keys) {
    $rowUpdates = $updates[$rowNum]
    foreach ($col in $rowUpdates.Keys) {
        $cellRef = "$col$rowNum"
        $cell = $ws.Range($cellRef)
        # Force text storage so values like "0.610" or "157.50" are not
        # reinterpreted as numbers (which would drop trailing zeros / change
        # formatting). Reset the style afterwards so no stray number-format
        # style lingers on the cell.
        $cell.NumberFormat = "@"
        $cell.Value = $rowUpdates[$col]
        $cell.Style = "Normal"
    }
}
